# Applies the "add gitignore and agg tab" structural change to the
# three worksheets of the workbook:
#   - inserts a new column I ("Mae pc от old ") that duplicates the
#     numeric Mae value currently held in column H
#   - renames the old H1 header from "Mae" to "Mae old от pc"
#   - the former column I ("Тип данных") is shifted one slot to the
#     right, to column J, automatically by the column insert
#   - on the "Negative Correlation" sheet the 19 data rows are
#     additionally reversed in place (A2:J20 -> A20:J2)

$wb = $excel.ActiveWorkbook

function Restructure-Sheet($ws) {
    $lastRow = $ws.UsedRange.Rows.Count

    # Insert a brand-new, empty column at I; this shifts the existing
    # column I ("Тип данных") one slot to the right, to J, together
    # with all of its header/body text, automatically.
    $ws.Columns.Item(9).Insert()

    # Update the two header cells in row 1 (J1 already holds the
    # "Тип данных" label that used to live in I1, courtesy of the
    # column insert above).
    $ws.Range("H1").Value = "Mae old от pc"
    $ws.Range("I1").Value = "Mae pc от old "

    # Copy H1's cell formatting (style) onto the freshly created I1
    # cell so it matches the rest of the header row.
    $ws.Range("H1").Copy()
    $ws.Range("I1").PasteSpecial(-4122)

    # Duplicate the numeric Mae values from column H into the new
    # column I for every data row.
    $hValues = $ws.Range("H2:H" + $lastRow).Value()
    $ws.Range("I2:I" + $lastRow).Value = $hValues

    return $lastRow
}

function Reverse-DataRows($ws, $firstRow, $lastRow, $lastColLetter) {
    $rangeAddr = "A" + $firstRow + ":" + $lastColLetter + $lastRow
    $arr = $ws.Range($rangeAddr).Value()

    $rLo = $arr.GetLowerBound(0)
    $rHi = $arr.GetUpperBound(0)
    $cLo = $arr.GetLowerBound(1)
    $cHi = $arr.GetUpperBound(1)
    $nRows = $rHi - $rLo + 1
    $nCols = $cHi - $cLo + 1

    $newArr = New-Object 'object[,]' $nRows, $nCols
    $rUpper = $newArr.GetUpperBound(0)
    $cUpper = $newArr.GetUpperBound(1)
    for ($r = 0; $r -le $rUpper; $r++) {
        $srcR = $rHi - $r
        for ($c = 0; $c -le $cUpper; $c++) {
            $srcC = $c + $cLo
            $val = $arr[$srcR, $srcC]
            $newArr[$r, $c] = $val
        }
    }
    $ws.Range($rangeAddr).Value = $newArr
}

$wsPositive = $wb.Worksheets.Item("Positive Correlation")
Restructure-Sheet $wsPositive | Out-Null

$wsNegative = $wb.Worksheets.Item("Negative Correlation")
$negLastRow = Restructure-Sheet $wsNegative
Reverse-DataRows $wsNegative 2 $negLastRow "J"

$wsGeneral = $wb.Worksheets.Item("General")
Restructure-Sheet $wsGeneral | Out-Null
